$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the UID value in A2 (kept the same shared-string slot, new text)
$ws.Range("A2").Value = """0000003"""

# Add a new "Note" column header in D1
$ws.Range("D1").Value = "Note"

# Leave D2 empty (no new data typed there yet) and move the selection onto it,
# matching the active cell left selected after the edit
$ws.Range("D2").Select()
